$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2809.8333
$ws.Range("J17").Value = 3035
$ws.Range("L17").Value = 9105
$ws.Range("N17").Value = -9441
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H58").Value = 254.75
$ws.Range("I58").Value = 73
$ws.Range("J58").Value = 800
$ws.Range("K58").Value = 219
$ws.Range("L58").Value = 2400
$ws.Range("M58").Value = -69
$ws.Range("N58").Value = -2700
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H96").Value = 884.2308
$ws.Range("J96").Value = 1562.6
$ws.Range("L96").Value = 4687.799999999999
$ws.Range("N96").Value = -7433.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 8000
$ws.Range("J106").Value = 8000
$ws.Range("L106").Value = 8000
$ws.Range("N106").Value = -10524
$ws.Range("H132").Value = 1656
$ws.Range("I132").Value = 1656
$ws.Range("K132").Value = 4968
$ws.Range("M132").Value = -2438

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 6250885
$ws.Range("I5").Value = 7143154
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 7143154
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = -7143041
$ws.Range("N5").Value = -5226
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H137").Value = 75000
$ws.Range("I137").Value = 50000
$ws.Range("J137").Value = 100000
$ws.Range("K137").Value = 50000
$ws.Range("L137").Value = 100000
$ws.Range("M137").Value = -44900
$ws.Range("N137").Value = -110200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2352.75
$ws.Range("I16").Value = 1498
$ws.Range("J16").Value = 2637.6667
$ws.Range("K16").Value = 1498
$ws.Range("L16").Value = 2637.6667
$ws.Range("M16").Value = -1211
$ws.Range("N16").Value = -3211.6667
$ws.Range("H28").Value = 22500
$ws.Range("J28").Value = 22500
$ws.Range("L28").Value = 22500
$ws.Range("N28").Value = -22990
$ws.Range("H113").Value = 2352.75
$ws.Range("I113").Value = 1498
$ws.Range("J113").Value = 2637.6667
$ws.Range("K113").Value = 1498
$ws.Range("L113").Value = 2637.6667
$ws.Range("M113").Value = 672
$ws.Range("N113").Value = -6977.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1582.3667
$ws.Range("I4").Value = 1228.375
$ws.Range("K4").Value = 3685.125
$ws.Range("M4").Value = -3573.125
$ws.Range("H10").Value = 62.782608
$ws.Range("I10").Value = 21.190475
$ws.Range("K10").Value = 63.571425
$ws.Range("M10").Value = 75.428575
$ws.Range("H23").Value = 499
$ws.Range("J23").Value = 498.33334
$ws.Range("L23").Value = 1495.00002
$ws.Range("N23").Value = -1965.00002
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H99").Value = 9258.166999999999
$ws.Range("I99").Value = 7637.5
$ws.Range("J99").Value = 12499.5
$ws.Range("K99").Value = 22912.5
$ws.Range("L99").Value = 37498.5
$ws.Range("M99").Value = -20666.5
$ws.Range("N99").Value = -41990.5
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("M117").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 31250
$ws.Range("J95").Value = 31250
$ws.Range("L95").Value = 31250
$ws.Range("N95").Value = -36742
$ws.Range("H102").Value = 2055.3
$ws.Range("I102").Value = 1298
$ws.Range("K102").Value = 1298
$ws.Range("M102").Value = 324

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 895
$ws.Range("I7").Value = 895
$ws.Range("K7").Value = 895
$ws.Range("M7").Value = -783
$ws.Range("H16").Value = 2000
$ws.Range("J16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("N16").Value = -2340
$ws.Range("H22").Value = 618.75
$ws.Range("I22").Value = 662.5
$ws.Range("J22").Value = 575
$ws.Range("K22").Value = 662.5
$ws.Range("L22").Value = 575
$ws.Range("M22").Value = -367.5
$ws.Range("N22").Value = -1165
$ws.Range("H27").Value = 618.75
$ws.Range("I27").Value = 662.5
$ws.Range("J27").Value = 575
$ws.Range("K27").Value = 662.5
$ws.Range("L27").Value = 575
$ws.Range("M27").Value = -555.5
$ws.Range("N27").Value = -789
$ws.Range("H35").Value = 12213.083
$ws.Range("I35").Value = 6093.6665
$ws.Range("K35").Value = 6093.6665
$ws.Range("M35").Value = -5757.6665
$ws.Range("H39").Value = 4933.3335
$ws.Range("I39").Value = 4933.3335
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 4933.3335
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -4473.3335
$ws.Range("N39").ClearContents()
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()
$ws.Range("H76").Value = 49739
$ws.Range("J76").Value = 49739
$ws.Range("L76").Value = 49739
$ws.Range("N76").Value = -50415
$ws.Range("H79").Value = 49739
$ws.Range("J79").Value = 49739
$ws.Range("L79").Value = 49739
$ws.Range("N79").Value = -52079
$ws.Range("H103").Value = 23999.5
$ws.Range("J103").Value = 23999.5
$ws.Range("L103").Value = 23999.5
$ws.Range("N103").Value = -26343.5
$ws.Range("H106").Value = 40999.668
$ws.Range("J106").Value = 40999.668
$ws.Range("L106").Value = 40999.668
$ws.Range("N106").Value = -43523.668
$ws.Range("H126").Value = 895
$ws.Range("I126").Value = 895
$ws.Range("K126").Value = 2685
$ws.Range("M126").Value = -215

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H95").Value = 8200
$ws.Range("J95").Value = 8200
$ws.Range("L95").Value = 8200
$ws.Range("N95").Value = -13692
$ws.Range("H113").Value = 862.5
$ws.Range("I113").Value = 461.5
$ws.Range("J113").Value = 1263.5
$ws.Range("K113").Value = 1384.5
$ws.Range("L113").Value = 3790.5
$ws.Range("M113").Value = 785.5
$ws.Range("N113").Value = -8130.5
$ws.Range("H117").Value = 42000
$ws.Range("J117").Value = 42000
$ws.Range("L117").Value = 42000
$ws.Range("N117").Value = -51178
